$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushing "last group"/"P" and subsequent rows down).
# Excel's row insert naturally carries the formatting of the row above into the new row.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(2).RowHeight

# Fill the new row 3 with "mid group" / "N"
$ws.Cells.Item(3, 1).Value = "mid group"
$ws.Cells.Item(3, 2).Value = "N"

# "questions per category"'s value cell (now row 7, column B) gets a new style:
# same font as before, plus left horizontal alignment.
$ws.Cells.Item(7, 2).HorizontalAlignment = -4131  # xlLeft

# Update the selection to B2 as in the diff
$ws.Range("B2").Select() | Out-Null
